$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2..96 hold labels q1..q95 in column A; shift each label down by one
# index (q1 -> q0, q2 -> q1, ..., q95 -> q94). Row 97 (q96) is left as-is.
for ($row = 2; $row -le 96; $row++) {
    $newIndex = $row - 2
    $ws.Cells.Item($row, 1).Value = "q$newIndex"
}
